$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'FAPs'
$ws.Cells.Item(2,2).Value = 'Fgf2'
$ws.Cells.Item(2,3).Value = 'Gpc4'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 13.41987666666667
$ws.Cells.Item(2,8).Value = 40.25962999999999
$ws.Cells.Item(2,9).Value = 0.744640910590638
$ws.Cells.Item(2,10).Value = 0.789423007068499
$ws.Cells.Item(2,11).Value = 2.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 17.853143
$ws.Cells.Item(2,14).Value = 35.706286
$ws.Cells.Item(2,15).Value = 0.263149724813253
$ws.Cells.Item(2,16).Value = 0.2070937143112604
$ws.Cells.Item(2,17).Value = 239.5869771723633
$ws.Cells.Item(2,18).Value = 1437.52186303418
$ws.Cells.Item(2,19).Value = 0.1959520507066165
$ws.Cells.Item(2,20).Value = 0.1634845426965799

# Row 3
$ws.Cells.Item(3,1).Value = 'FAPs'
$ws.Cells.Item(3,2).Value = 'Fgf2'
$ws.Cells.Item(3,3).Value = 'Gpc4'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 13.41987666666667
$ws.Cells.Item(3,8).Value = 40.25962999999999
$ws.Cells.Item(3,9).Value = 0.744640910590638
$ws.Cells.Item(3,10).Value = 0.789423007068499
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 25.610932
$ws.Cells.Item(3,14).Value = 76.832796
$ws.Cells.Item(3,15).Value = 0.3774970999790309
$ws.Cells.Item(3,16).Value = 0.4456243112083781
$ws.Cells.Item(3,17).Value = 343.6955487583867
$ws.Cells.Item(3,18).Value = 3093.25993882548
$ws.Cells.Item(3,19).Value = 0.2810997842737106
$ws.Cells.Item(3,20).Value = 0.3517860837769465

# Row 4
$ws.Cells.Item(4,1).Value = 'FAPs'
$ws.Cells.Item(4,2).Value = 'Fgf2'
$ws.Cells.Item(4,3).Value = 'Gpc4'
$ws.Cells.Item(4,4).Value = 'Neutro'
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 13.41987666666667
$ws.Cells.Item(4,8).Value = 40.25962999999999
$ws.Cells.Item(4,9).Value = 0.744640910590638
$ws.Cells.Item(4,10).Value = 0.789423007068499
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 11.11704666666667
$ws.Cells.Item(4,14).Value = 33.35114
$ws.Cells.Item(4,15).Value = 0.1638617789074688
$ws.Cells.Item(4,16).Value = 0.1934340485346152
$ws.Cells.Item(4,17).Value = 149.1893951642444
$ws.Cells.Item(4,18).Value = 1342.7045564782
$ws.Cells.Item(4,19).Value = 0.1220181842566594
$ws.Cells.Item(4,20).Value = 0.1527012882636299

# Row 5
$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'Fgf2'
$ws.Cells.Item(5,3).Value = 'Gpc4'
$ws.Cells.Item(5,4).Value = 'sCs'
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 13.41987666666667
$ws.Cells.Item(5,8).Value = 40.25962999999999
$ws.Cells.Item(5,9).Value = 0.744640910590638
$ws.Cells.Item(5,10).Value = 0.789423007068499
$ws.Cells.Item(5,11).Value = 2.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 13.262928
$ws.Cells.Item(5,14).Value = 26.525856
$ws.Cells.Item(5,15).Value = 0.1954913963002474
$ws.Cells.Item(5,16).Value = 0.1538479259457462
$ws.Cells.Item(5,17).Value = 177.98685799888
$ws.Cells.Item(5,18).Value = 1067.92114799328
$ws.Cells.Item(5,19).Value = 0.1455708913536515
$ws.Cells.Item(5,20).Value = 0.1214510923313427

# Row 6
$ws.Cells.Item(6,1).Value = 'M1'
$ws.Cells.Item(6,2).Value = 'Fgf2'
$ws.Cells.Item(6,3).Value = 'Gpc4'
$ws.Cells.Item(6,4).Value = 'ECs'
$ws.Cells.Item(6,5).Value = 1.0
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.1169053333333333
$ws.Cells.Item(6,8).Value = 0.350716
$ws.Cells.Item(6,9).Value = 0.006486832631067555
$ws.Cells.Item(6,10).Value = 0.00687694544999633
$ws.Cells.Item(6,11).Value = 2.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 17.853143
$ws.Cells.Item(6,14).Value = 35.706286
$ws.Cells.Item(6,15).Value = 0.263149724813253
$ws.Cells.Item(6,16).Value = 0.2070937143112604
$ws.Cells.Item(6,17).Value = 2.087127633462667
$ws.Cells.Item(6,18).Value = 12.522765800776
$ws.Cells.Item(6,19).Value = 0.001707008221775057
$ws.Cells.Item(6,20).Value = 0.001424172176355662

# Row 7
$ws.Cells.Item(7,1).Value = 'M1'
$ws.Cells.Item(7,2).Value = 'Fgf2'
$ws.Cells.Item(7,3).Value = 'Gpc4'
$ws.Cells.Item(7,4).Value = 'FAPs'
$ws.Cells.Item(7,5).Value = 1.0
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.1169053333333333
$ws.Cells.Item(7,8).Value = 0.350716
$ws.Cells.Item(7,9).Value = 0.006486832631067555
$ws.Cells.Item(7,10).Value = 0.00687694544999633
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 25.610932
$ws.Cells.Item(7,14).Value = 76.832796
$ws.Cells.Item(7,15).Value = 0.3774970999790309
$ws.Cells.Item(7,16).Value = 0.4456243112083781
$ws.Cells.Item(7,17).Value = 2.994054542437334
$ws.Cells.Item(7,18).Value = 26.946490881936
$ws.Cells.Item(7,19).Value = 0.002448760506277349
$ws.Cells.Item(7,20).Value = 0.003064534079372204

# Row 8
$ws.Cells.Item(8,1).Value = 'M1'
$ws.Cells.Item(8,2).Value = 'Fgf2'
$ws.Cells.Item(8,3).Value = 'Gpc4'
$ws.Cells.Item(8,4).Value = 'Neutro'
$ws.Cells.Item(8,5).Value = 1.0
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.1169053333333333
$ws.Cells.Item(8,8).Value = 0.350716
$ws.Cells.Item(8,9).Value = 0.006486832631067555
$ws.Cells.Item(8,10).Value = 0.00687694544999633
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 11.11704666666667
$ws.Cells.Item(8,14).Value = 33.35114
$ws.Cells.Item(8,15).Value = 0.1638617789074688
$ws.Cells.Item(8,16).Value = 0.1934340485346152
$ws.Cells.Item(8,17).Value = 1.299642046248889
$ws.Cells.Item(8,18).Value = 11.69677841624
$ws.Cells.Item(8,19).Value = 0.001062943934401746
$ws.Cells.Item(8,20).Value = 0.001330235399944491

# Row 9
$ws.Cells.Item(9,1).Value = 'M1'
$ws.Cells.Item(9,2).Value = 'Fgf2'
$ws.Cells.Item(9,3).Value = 'Gpc4'
$ws.Cells.Item(9,4).Value = 'sCs'
$ws.Cells.Item(9,5).Value = 1.0
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.1169053333333333
$ws.Cells.Item(9,8).Value = 0.350716
$ws.Cells.Item(9,9).Value = 0.006486832631067555
$ws.Cells.Item(9,10).Value = 0.00687694544999633
$ws.Cells.Item(9,11).Value = 2.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 13.262928
$ws.Cells.Item(9,14).Value = 26.525856
$ws.Cells.Item(9,15).Value = 0.1954913963002474
$ws.Cells.Item(9,16).Value = 0.1538479259457462
$ws.Cells.Item(9,17).Value = 1.550507018816
$ws.Cells.Item(9,18).Value = 9.303042112896001
$ws.Cells.Item(9,19).Value = 0.001268119968613404
$ws.Cells.Item(9,20).Value = 0.001058003794323972

# Row 10
$ws.Cells.Item(10,1).Value = 'M2'
$ws.Cells.Item(10,2).Value = 'Fgf2'
$ws.Cells.Item(10,3).Value = 'Gpc4'
$ws.Cells.Item(10,4).Value = 'ECs'
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 0.5570086666666666
$ws.Cells.Item(10,8).Value = 1.671026
$ws.Cells.Item(10,9).Value = 0.03090724684406269
$ws.Cells.Item(10,10).Value = 0.03276598343823939
$ws.Cells.Item(10,11).Value = 2.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 17.853143
$ws.Cells.Item(10,14).Value = 35.706286
$ws.Cells.Item(10,15).Value = 0.263149724813253
$ws.Cells.Item(10,16).Value = 0.2070937143112604
$ws.Cells.Item(10,17).Value = 9.944355378239331
$ws.Cells.Item(10,18).Value = 59.66613226943599
$ws.Cells.Item(10,19).Value = 0.008133233501750377
$ws.Cells.Item(10,20).Value = 0.006785629213286239

# Row 11
$ws.Cells.Item(11,1).Value = 'M2'
$ws.Cells.Item(11,2).Value = 'Fgf2'
$ws.Cells.Item(11,3).Value = 'Gpc4'
$ws.Cells.Item(11,4).Value = 'FAPs'
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 0.5570086666666666
$ws.Cells.Item(11,8).Value = 1.671026
$ws.Cells.Item(11,9).Value = 0.03090724684406269
$ws.Cells.Item(11,10).Value = 0.03276598343823939
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 25.610932
$ws.Cells.Item(11,14).Value = 76.832796
$ws.Cells.Item(11,15).Value = 0.3774970999790309
$ws.Cells.Item(11,16).Value = 0.4456243112083781
$ws.Cells.Item(11,17).Value = 14.26551108541067
$ws.Cells.Item(11,18).Value = 128.389599768696
$ws.Cells.Item(11,19).Value = 0.01166739605196972
$ws.Cells.Item(11,20).Value = 0.01460131880073055

# Row 12
$ws.Cells.Item(12,1).Value = 'M2'
$ws.Cells.Item(12,2).Value = 'Fgf2'
$ws.Cells.Item(12,3).Value = 'Gpc4'
$ws.Cells.Item(12,4).Value = 'Neutro'
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 0.5570086666666666
$ws.Cells.Item(12,8).Value = 1.671026
$ws.Cells.Item(12,9).Value = 0.03090724684406269
$ws.Cells.Item(12,10).Value = 0.03276598343823939
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 11.11704666666667
$ws.Cells.Item(12,14).Value = 33.35114
$ws.Cells.Item(12,15).Value = 0.1638617789074688
$ws.Cells.Item(12,16).Value = 0.1934340485346152
$ws.Cells.Item(12,17).Value = 6.192291341071111
$ws.Cells.Item(12,18).Value = 55.73062206964
$ws.Cells.Item(12,19).Value = 0.005064516449000364
$ws.Cells.Item(12,20).Value = 0.006338056830676795

# Row 13
$ws.Cells.Item(13,1).Value = 'M2'
$ws.Cells.Item(13,2).Value = 'Fgf2'
$ws.Cells.Item(13,3).Value = 'Gpc4'
$ws.Cells.Item(13,4).Value = 'sCs'
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 0.5570086666666666
$ws.Cells.Item(13,8).Value = 1.671026
$ws.Cells.Item(13,9).Value = 0.03090724684406269
$ws.Cells.Item(13,10).Value = 0.03276598343823939
$ws.Cells.Item(13,11).Value = 2.0
$ws.Cells.Item(13,12).Value = 1.0
$ws.Cells.Item(13,13).Value = 13.262928
$ws.Cells.Item(13,14).Value = 26.525856
$ws.Cells.Item(13,15).Value = 0.1954913963002474
$ws.Cells.Item(13,16).Value = 0.1538479259457462
$ws.Cells.Item(13,17).Value = 7.387565841375999
$ws.Cells.Item(13,18).Value = 44.325395048256
$ws.Cells.Item(13,19).Value = 0.006042100841342229
$ws.Cells.Item(13,20).Value = 0.0050409785935458

# Row 14
$ws.Cells.Item(14,1).Value = 'Neutro'
$ws.Cells.Item(14,2).Value = 'Fgf2'
$ws.Cells.Item(14,3).Value = 'Gpc4'
$ws.Cells.Item(14,4).Value = 'ECs'
$ws.Cells.Item(14,5).Value = 3.0
$ws.Cells.Item(14,6).Value = 1.0
$ws.Cells.Item(14,7).Value = 0.8611266666666667
$ws.Cells.Item(14,8).Value = 2.58338
$ws.Cells.Item(14,9).Value = 0.04778211910048957
$ws.Cells.Item(14,10).Value = 0.05065569673642354
$ws.Cells.Item(14,11).Value = 2.0
$ws.Cells.Item(14,12).Value = 1.0
$ws.Cells.Item(14,13).Value = 17.853143
$ws.Cells.Item(14,14).Value = 35.706286
$ws.Cells.Item(14,15).Value = 0.263149724813253
$ws.Cells.Item(14,16).Value = 0.2070937143112604
$ws.Cells.Item(14,17).Value = 15.37381752111333
$ws.Cells.Item(14,18).Value = 92.24290512668
$ws.Cells.Item(14,19).Value = 0.01257385149228791
$ws.Cells.Item(14,20).Value = 0.01049047638817074

# Row 15
$ws.Cells.Item(15,1).Value = 'Neutro'
$ws.Cells.Item(15,2).Value = 'Fgf2'
$ws.Cells.Item(15,3).Value = 'Gpc4'
$ws.Cells.Item(15,4).Value = 'FAPs'
$ws.Cells.Item(15,5).Value = 3.0
$ws.Cells.Item(15,6).Value = 1.0
$ws.Cells.Item(15,7).Value = 0.8611266666666667
$ws.Cells.Item(15,8).Value = 2.58338
$ws.Cells.Item(15,9).Value = 0.04778211910048957
$ws.Cells.Item(15,10).Value = 0.05065569673642354
$ws.Cells.Item(15,11).Value = 3.0
$ws.Cells.Item(15,12).Value = 1.0
$ws.Cells.Item(15,13).Value = 25.610932
$ws.Cells.Item(15,14).Value = 76.832796
$ws.Cells.Item(15,15).Value = 0.3774970999790309
$ws.Cells.Item(15,16).Value = 0.4456243112083781
$ws.Cells.Item(15,17).Value = 22.05425650338667
$ws.Cells.Item(15,18).Value = 198.48830853048
$ws.Cells.Item(15,19).Value = 0.01803761139128747
$ws.Cells.Item(15,20).Value = 0.02257340996694923

# Row 16
$ws.Cells.Item(16,1).Value = 'Neutro'
$ws.Cells.Item(16,2).Value = 'Fgf2'
$ws.Cells.Item(16,3).Value = 'Gpc4'
$ws.Cells.Item(16,4).Value = 'Neutro'
$ws.Cells.Item(16,5).Value = 3.0
$ws.Cells.Item(16,6).Value = 1.0
$ws.Cells.Item(16,7).Value = 0.8611266666666667
$ws.Cells.Item(16,8).Value = 2.58338
$ws.Cells.Item(16,9).Value = 0.04778211910048957
$ws.Cells.Item(16,10).Value = 0.05065569673642354
$ws.Cells.Item(16,11).Value = 3.0
$ws.Cells.Item(16,12).Value = 1.0
$ws.Cells.Item(16,13).Value = 11.11704666666667
$ws.Cells.Item(16,14).Value = 33.35114
$ws.Cells.Item(16,15).Value = 0.1638617789074688
$ws.Cells.Item(16,16).Value = 0.1934340485346152
$ws.Cells.Item(16,17).Value = 9.573185339244445
$ws.Cells.Item(16,18).Value = 86.1586680532
$ws.Cells.Item(16,19).Value = 0.007829663035774766
$ws.Cells.Item(16,20).Value = 0.0097985365010681

# Row 17
$ws.Cells.Item(17,1).Value = 'Neutro'
$ws.Cells.Item(17,2).Value = 'Fgf2'
$ws.Cells.Item(17,3).Value = 'Gpc4'
$ws.Cells.Item(17,4).Value = 'sCs'
$ws.Cells.Item(17,5).Value = 3.0
$ws.Cells.Item(17,6).Value = 1.0
$ws.Cells.Item(17,7).Value = 0.8611266666666667
$ws.Cells.Item(17,8).Value = 2.58338
$ws.Cells.Item(17,9).Value = 0.04778211910048957
$ws.Cells.Item(17,10).Value = 0.05065569673642354
$ws.Cells.Item(17,11).Value = 2.0
$ws.Cells.Item(17,12).Value = 1.0
$ws.Cells.Item(17,13).Value = 13.262928
$ws.Cells.Item(17,14).Value = 26.525856
$ws.Cells.Item(17,15).Value = 0.1954913963002474
$ws.Cells.Item(17,16).Value = 0.1538479259457462
$ws.Cells.Item(17,17).Value = 11.42106097888
$ws.Cells.Item(17,18).Value = 68.52636587328
$ws.Cells.Item(17,19).Value = 0.009340993181139427
$ws.Cells.Item(17,20).Value = 0.007793273880235466

# Row 18
$ws.Cells.Item(18,1).Value = 'sCs'
$ws.Cells.Item(18,2).Value = 'Fgf2'
$ws.Cells.Item(18,3).Value = 'Gpc4'
$ws.Cells.Item(18,4).Value = 'ECs'
$ws.Cells.Item(18,5).Value = 2.0
$ws.Cells.Item(18,6).Value = 1.0
$ws.Cells.Item(18,7).Value = 3.0670265
$ws.Cells.Item(18,8).Value = 6.134053
$ws.Cells.Item(18,9).Value = 0.1701828908337422
$ws.Cells.Item(18,10).Value = 0.1202783673068418
$ws.Cells.Item(18,11).Value = 2.0
$ws.Cells.Item(18,12).Value = 1.0
$ws.Cells.Item(18,13).Value = 17.853143
$ws.Cells.Item(18,14).Value = 35.706286
$ws.Cells.Item(18,15).Value = 0.263149724813253
$ws.Cells.Item(18,16).Value = 0.2070937143112604
$ws.Cells.Item(18,17).Value = 54.75606268928949
$ws.Cells.Item(18,18).Value = 219.024250757158
$ws.Cells.Item(18,19).Value = 0.04478358089082314
$ws.Cells.Item(18,20).Value = 0.02490889383686795

# Row 19
$ws.Cells.Item(19,1).Value = 'sCs'
$ws.Cells.Item(19,2).Value = 'Fgf2'
$ws.Cells.Item(19,3).Value = 'Gpc4'
$ws.Cells.Item(19,4).Value = 'FAPs'
$ws.Cells.Item(19,5).Value = 2.0
$ws.Cells.Item(19,6).Value = 1.0
$ws.Cells.Item(19,7).Value = 3.0670265
$ws.Cells.Item(19,8).Value = 6.134053
$ws.Cells.Item(19,9).Value = 0.1701828908337422
$ws.Cells.Item(19,10).Value = 0.1202783673068418
$ws.Cells.Item(19,11).Value = 3.0
$ws.Cells.Item(19,12).Value = 1.0
$ws.Cells.Item(19,13).Value = 25.610932
$ws.Cells.Item(19,14).Value = 76.832796
$ws.Cells.Item(19,15).Value = 0.3774970999790309
$ws.Cells.Item(19,16).Value = 0.4456243112083781
$ws.Cells.Item(19,17).Value = 78.549407133698
$ws.Cells.Item(19,18).Value = 471.296442802188
$ws.Cells.Item(19,19).Value = 0.06424354775578568
$ws.Cells.Item(19,20).Value = 0.05359896458437969

# Row 20
$ws.Cells.Item(20,1).Value = 'sCs'
$ws.Cells.Item(20,2).Value = 'Fgf2'
$ws.Cells.Item(20,3).Value = 'Gpc4'
$ws.Cells.Item(20,4).Value = 'Neutro'
$ws.Cells.Item(20,5).Value = 2.0
$ws.Cells.Item(20,6).Value = 1.0
$ws.Cells.Item(20,7).Value = 3.0670265
$ws.Cells.Item(20,8).Value = 6.134053
$ws.Cells.Item(20,9).Value = 0.1701828908337422
$ws.Cells.Item(20,10).Value = 0.1202783673068418
$ws.Cells.Item(20,11).Value = 3.0
$ws.Cells.Item(20,12).Value = 1.0
$ws.Cells.Item(20,13).Value = 11.11704666666667
$ws.Cells.Item(20,14).Value = 33.35114
$ws.Cells.Item(20,15).Value = 0.1638617789074688
$ws.Cells.Item(20,16).Value = 0.1934340485346152
$ws.Cells.Item(20,17).Value = 34.09627672840333
$ws.Cells.Item(20,18).Value = 204.57766037042
$ws.Cells.Item(20,19).Value = 0.02788647123163257
$ws.Cells.Item(20,20).Value = 0.02326593153929591

# Row 21
$ws.Cells.Item(21,1).Value = 'sCs'
$ws.Cells.Item(21,2).Value = 'Fgf2'
$ws.Cells.Item(21,3).Value = 'Gpc4'
$ws.Cells.Item(21,4).Value = 'sCs'
$ws.Cells.Item(21,5).Value = 2.0
$ws.Cells.Item(21,6).Value = 1.0
$ws.Cells.Item(21,7).Value = 3.0670265
$ws.Cells.Item(21,8).Value = 6.134053
$ws.Cells.Item(21,9).Value = 0.1701828908337422
$ws.Cells.Item(21,10).Value = 0.1202783673068418
$ws.Cells.Item(21,11).Value = 2.0
$ws.Cells.Item(21,12).Value = 1.0
$ws.Cells.Item(21,13).Value = 13.262928
$ws.Cells.Item(21,14).Value = 26.525856
$ws.Cells.Item(21,15).Value = 0.1954913963002474
$ws.Cells.Item(21,16).Value = 0.1538479259457462
$ws.Cells.Item(21,17).Value = 40.677751643592
$ws.Cells.Item(21,18).Value = 162.711006574368
$ws.Cells.Item(21,19).Value = 0.03326929095550084
$ws.Cells.Item(21,20).Value = 0.01850457734629826
